# "Generate Report for Handoff"
#
# The localization-status report was regenerated: the two files that were
# still "Ready for handoff" (52e2b374.../b659bb28.../c81de128.../cea53375...)
# had their priority bumped from "low" to "ht", and their handoff/generate
# timestamps refreshed.

$wb = $excel.ActiveWorkbook

$ovr   = $wb.Worksheets.Item("Overview")
$zhcn  = $wb.Worksheets.Item("zh-cn")
$dede  = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" (column G) refreshed for rows 4-7
$ovr.Range("G4").Value = "2016-09-02 14:37:57"
$ovr.Range("G5").Value = "2016-09-02 14:37:57"
$ovr.Range("G6").Value = "2016-09-02 14:37:57"
$ovr.Range("G7").Value = "2016-09-02 14:37:57"

# zh-cn sheet: Priority (E) low -> ht, Latest Handoff Datetime (H) refreshed
$zhcn.Range("E4").Value = "ht"
$zhcn.Range("E5").Value = "ht"
$zhcn.Range("E6").Value = "ht"
$zhcn.Range("E7").Value = "ht"

$zhcn.Range("H4").Value = "2016-09-02 14:37:52"
$zhcn.Range("H5").Value = "2016-09-02 14:37:52"
$zhcn.Range("H6").Value = "2016-09-02 14:37:52"
$zhcn.Range("H7").Value = "2016-09-02 14:37:52"

# de-de sheet: Priority (E) low -> ht
$dede.Range("E4").Value = "ht"
$dede.Range("E5").Value = "ht"
$dede.Range("E6").Value = "ht"
$dede.Range("E7").Value = "ht"

# de-de sheet shares the same "Latest HO Xliff Generate Date" text as the
# Overview sheet (column H here) - refresh it too so the old timestamp
# string is fully superseded (matches how the shared-string table collapses
# in the real commit).
$dede.Range("H4").Value = "2016-09-02 14:37:57"
$dede.Range("H5").Value = "2016-09-02 14:37:57"
$dede.Range("H6").Value = "2016-09-02 14:37:57"
$dede.Range("H7").Value = "2016-09-02 14:37:57"
